# Commit: "zh-tW -> Eng to test embedding, and refind qa_type"
# Translate the qa_type column (column B) from Traditional Chinese to English:
#   問答 (qa/Q&A)            -> qa
#   是非 (true/false)        -> true_false
#   選擇 (multiple choice)   -> multiple_choice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "問答" = "qa"
    "是非" = "true_false"
    "選擇" = "multiple_choice"
}

# Data rows run from row 2 through row 56 (row 1 is the header: id/type/question/answers/keywords)
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# Update the window view: scroll so row 15 is at the top and select B20
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("B20").Select()
